$d = $word.ActiveDocument

# Bab V, point 2: restructure the sentence about the SG90 servo.
# Old: "Penelitian yang dilakukan menggunakan servo SG90 yang memiliki gir ..."
# New: "Penggunaan servo SG90 dalam penelitian yang dilakukan memiliki gir ..."

$old = "Penelitian yang dilakukan menggunakan servo SG90 yang memiliki gir"
$new = "Penggunaan servo SG90 dalam penelitian yang dilakukan memiliki gir"

$found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Target sentence for Bab V point 2 was not found; aborting to avoid silent no-op."
}

Write-Host "Replaced: $found"
